$wb = $excel.ActiveWorkbook

# --- Sheet1: update the maxrows/noheader-maxrows template strings to include
# the "Other" aggregation bucket + fill value, and collapse the "noheader"
# block to a single-style block like the real data rows (no separate header
# row) since no-header dfs no longer render their own header styling.
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Range("A4").Value = "{{ df2 | maxrows(2, Other, 0)}}"
$ws1.Range("A7").Value = "{{ df2 | noheader | maxrows(2, Other, 0) }}"

# Row 7 (the "noheader" frame row) now uses the plain data-row style (style 4)
# instead of the header style (style 3) - copy formatting from row 5.
$ws1.Range("A5:E5").Copy()
$ws1.Range("A7:E7").PasteSpecial(-4122)

# The noheader block no longer reserves a separate row for its data (row 8) -
# it renders directly under row 7, so remove the now-unused row 8.
$ws1.Rows.Item(8).Delete()

# --- Sheet2 ("expected"): the maxrows(2, Other, 0) aggregation means that
# once a dataframe has more rows than maxrows, the extra rows get collapsed
# into a single "Other" bucket row (replacing what used to be the literal
# index value 1).
$ws2 = $wb.Worksheets.Item("expected")

$ws2.Range("A5:E5").Copy()
$ws2.Range("A7:E7").PasteSpecial(-4122)

$ws2.Range("A5").Value = "Other"
$ws2.Range("A8").Value = "Other"

Write-Output "edits applied"
